$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 'FDEMO'
$ws.Range("B5").Value = '5DEMO'
$ws.Range("B6").Value = 'K312'
$ws.Range("B7").Value = 'P335'
$ws.Range("A9").Value = 'Fyziologie živočichů a člověka'
$ws.Range("B9").Value = 'PB502'
$ws.Range("A10").Value = 'Diplomová práce I'
$ws.Range("B10").Value = 'N028'
$ws.Range("A11").Value = 'Fyziologie živočichů a člověka'
$ws.Range("B11").Value = 'KB502'
$ws.Range("A12").Value = 'Diplomová práce I'
$ws.Range("B12").Value = 'M200'
$ws.Range("A13").Value = 'Proseminář řešení fyzikálních úloh A ZŠ'
$ws.Range("B13").Value = 'P736'
$ws.Range("A14").Value = 'Fyzikální praktikum C'
$ws.Range("B14").Value = 'K521'
$ws.Range("A15").Value = 'Didaktická a laboratorní technika'
$ws.Range("B15").Value = 'PA31'
$ws.Range("A16").Value = 'Souvislá pedagogická praxe II - ZŠ'
$ws.Range("B16").Value = 'KA35'
$ws.Range("A18").Value = 'Exkurze'
$ws.Range("B18").Value = 'P841'
$ws.Range("A19").Value = 'Proseminář řešení fyzikálních úloh'
$ws.Range("B19").Value = 'K622'
$ws.Range("A22").Value = 'Vybrané partie z fyziky'
$ws.Range("B22").Value = 'K414'
$ws.Range("A23").Value = 'Matematický proseminář'
$ws.Range("B23").Value = 'P365'
$ws.Range("A24").Value = 'Kvantová fyzika'
$ws.Range("B24").Value = 'P737'
$ws.Range("A26").Value = 'Vybrané partie z teorie elektromag. pole'
$ws.Range("B26").Value = 'PD203'
$ws.Range("A27").Value = 'Seminář katedry II'
$ws.Range("B27").Value = 'PA52'
$ws.Range("A28").Value = 'Seminář katedry'
$ws.Range("B28").Value = 'P602'
$ws.Range("A29").Value = 'Kvantová mechanika'
$ws.Range("B29").Value = 'KKM'
$ws.Range("A30").Value = 'Seminář katedry I'
$ws.Range("B30").Value = 'P952'
$ws.Range("A31").Value = 'Odborný seminář katedry A'
$ws.Range("B31").Value = 'P531'
$ws.Range("A32").Value = 'Kvantová fyzika I'
$ws.Range("B32").Value = 'K512'
$ws.Range("A33").Value = 'PVK - Vyb.partie z teoretické fyziky I'
$ws.Range("B33").Value = 'K422'
$ws.Range("A34").Value = 'PVK - Vyb.partie z teoretické fyziky II'
$ws.Range("B34").Value = 'K522'
$ws.Range("A35").Value = 'Teorie elektromagnetického pole'
$ws.Range("B35").Value = 'P508'
$ws.Range("A36").Value = 'Fyzika I'
$ws.Range("B36").Value = 'K222'
$ws.Range("A37").Value = 'Souhrnný seminář - počítačové modelování'
$ws.Range("B37").Value = 'SPM'
$ws.Range("A38").Value = 'Kvantová mechanika'
$ws.Range("B38").Value = 'KM'
$ws.Range("A39").Value = 'Teorie elektromagnetického pole'
$ws.Range("B39").Value = 'TEP'
$ws.Range("A40").Value = 'Vybrané partie z teoretické fyziky I'
$ws.Range("B40").Value = 'P938'
$ws.Range("A41").Value = 'Kvantová fyzika I'
$ws.Range("B41").Value = 'P509'
$ws.Range("A42").Value = 'Souhrnný seminář - počítačové modelování'
$ws.Range("B42").Value = 'KSPM'
$ws.Range("A43").Value = 'Seminář katedry'
$ws.Range("B43").Value = 'K602'
$ws.Range("A44").Value = 'Kvantová fyzika II'
$ws.Range("B44").Value = 'P607'
$ws.Range("B49").Value = 'N304'
$ws.Range("B51").Value = 'KN22'
$ws.Range("A53").Value = 'Vybrané matematické metody ve fyzice'
$ws.Range("B53").Value = '0212'
$ws.Range("A54").Value = 'VK- Praktikum mikrovlnné techniky D'
$ws.Range("B54").Value = '0211'
$ws.Range("A55").Value = 'Praktikum mikrovlnné techniky - aplikace'
$ws.Range("B55").Value = '0207'
$ws.Range("A56").Value = 'VK-Programování ATMEL I'
$ws.Range("B56").Value = '0175'
$ws.Range("A57").Value = 'Experimentální využití elmg. vln B'
$ws.Range("B57").Value = '0218'
$ws.Range("A58").Value = 'Experimentální využití elmg. vln F'
$ws.Range("B58").Value = '0237'
$ws.Range("A59").Value = 'Vybrané partie z teoretické fyziky I'
$ws.Range("B59").Value = 'P422'
$ws.Range("A60").Value = 'Vybrané partie z teoretické fyziky II'
$ws.Range("B60").Value = 'P522'
$ws.Range("B62").Value = 'P103'
$ws.Range("A63").Value = 'PVK - Hromadné zpracování dat'
$ws.Range("B63").Value = 'K406'
$ws.Range("A65").Value = 'Programování B'
$ws.Range("B65").Value = 'P203'
$ws.Range("A66").Value = 'Programování A'
$ws.Range("B66").Value = 'K103'
$ws.Range("B68").Value = 'P406'
$ws.Range("A69").Value = 'VK-Fyzikální seminář'
$ws.Range("B69").Value = '0204'
$ws.Range("A70").Value = 'Počítačové modelování II'
$ws.Range("B70").Value = 'K206'
$ws.Range("A71").Value = 'Počítačové modelování III'
$ws.Range("B71").Value = 'P401'
$ws.Range("A81").Value = 'Obecná zoologie'
$ws.Range("B81").Value = 'K110'
$ws.Range("A82").Value = 'Ekofyziologie'
$ws.Range("B82").Value = 'N051'
$ws.Range("A83").Value = 'Biologie a ekologie člověka I'
$ws.Range("B83").Value = 'P303'
$ws.Range("A84").Value = 'Somatologie'
$ws.Range("B84").Value = 'BP423'
$ws.Range("A85").Value = 'Biologie parazitů'
$ws.Range("B85").Value = 'N032'
$ws.Range("A86").Value = 'Obecná parazitologie'
$ws.Range("B86").Value = 'P323'
$ws.Range("A87").Value = 'Somatologie'
$ws.Range("B87").Value = 'BK423'
$ws.Range("A88").Value = 'Komplexní analýza materiálů'
$ws.Range("B88").Value = 'PD202'
$ws.Range("A89").Value = "Met. analýz mat. `nvyuž. el. mikroskopie"
$ws.Range("B89").Value = 'AP03'
$ws.Range("A91").Value = 'Kvantová fyzika I'
$ws.Range("B91").Value = 'P509'
$ws.Range("B92").Value = 'ME200'
$ws.Range("B93").Value = 'N023'
$ws.Range("A94").Value = 'Limnobiologie'
$ws.Range("B94").Value = 'P528'
$ws.Range("B95").Value = 'P328'
$ws.Range("A97").Value = 'Limnobiologie'
$ws.Range("B97").Value = 'BK108'
$ws.Range("A98").Value = 'Algologie'
$ws.Range("B98").Value = 'P420'
$ws.Range("B99").Value = 'P525'
$ws.Range("A100").Value = 'Environmentální legislativa'
$ws.Range("B100").Value = 'K507'
$ws.Range("B101").Value = 'N036'
$ws.Range("B102").Value = 'N001'
$ws.Range("B103").Value = 'M105'
$ws.Range("A108").Value = 'Souvislá pedagogická praxe I SŠ'
$ws.Range("B108").Value = 'P706'
$ws.Range("A109").Value = 'Biologie a ekologie člověka I'
$ws.Range("B109").Value = 'P303'
$ws.Range("A112").Value = 'VK - Konverzace němčiny pro geografy'
$ws.Range("B112").Value = '0130'
$ws.Range("A113").Value = 'VK-Němčina pro geografy'
$ws.Range("B113").Value = '0126'
$ws.Range("B114").Value = 'ME200'
$ws.Range("B115").Value = 'N023'
$ws.Range("B123").Value = 'P419'
$ws.Range("B124").Value = 'K403'
$ws.Range("B126").Value = 'P419'
$ws.Range("B127").Value = 'K403'
$ws.Range("A134").Value = 'Matematika I'
$ws.Range("B134").Value = 'P106'
$ws.Range("A135").Value = 'Úvod do fyziky'
$ws.Range("B135").Value = 'P112'
$ws.Range("A137").Value = 'Kvantová fyzika I'
$ws.Range("B137").Value = 'P509'
$ws.Range("A139").Value = 'Fyzikálně chemická cvičení'
$ws.Range("B139").Value = 'P527'
$ws.Range("A140").Value = 'Moderní elektroanal. met. - prakt. kurz'
$ws.Range("B140").Value = 'N029'
$ws.Range("A141").Value = 'PVK - Fyzikálně chemická cvičení A'
$ws.Range("B141").Value = 'P506'
$ws.Range("A142").Value = 'Laboratorní cvičení z fyzikální chemie'
$ws.Range("B142").Value = 'P532'
$ws.Range("A143").Value = 'Fyzikálně chemická cvičení'
$ws.Range("B143").Value = 'K502'
$ws.Range("A144").Value = 'Fyzikálně chemická cvičení'
$ws.Range("B144").Value = 'P517'
$ws.Range("B148").Value = 'N005'
$ws.Range("B149").Value = 'KN46'
$ws.Range("A151").Value = 'Mikroprocesory a senzory v praxi I'
$ws.Range("B151").Value = '0164'
$ws.Range("A152").Value = 'Mikroprocesory a senzory v praxi II'
$ws.Range("B152").Value = '0171'
$ws.Range("B161").Value = 'BK401'
$ws.Range("A164").Value = 'Životní prostředí a udržitelný rozvoj'
$ws.Range("B164").Value = 'B401'
$ws.Range("A165").Value = 'Metody geografického výzkumu krajiny'
$ws.Range("B165").Value = 'M203'
$ws.Range("A168").Value = 'Zoologické terénní cvičení'
$ws.Range("B168").Value = 'PB427'
$ws.Range("A169").Value = 'Ornitologie'
$ws.Range("B169").Value = 'MA205'
$ws.Range("B170").Value = 'KB427'
$ws.Range("A182").Value = 'Průběžná pedagogická praxe'
$ws.Range("B182").Value = 'P705'
$ws.Range("A183").Value = 'Didaktika biologie I'
$ws.Range("B183").Value = 'P108'
$ws.Range("A185").Value = 'Didaktika biologie II'
$ws.Range("B185").Value = 'P208'
$ws.Range("A186").Value = 'Souvislá pedagogická praxe II - SŠ'
$ws.Range("B186").Value = 'P707B'
$ws.Range("A187").Value = 'Souvislá pedagogická praxe I - SŠ'
$ws.Range("B187").Value = 'P707A'
$ws.Range("A190").Value = 'Apl. biologie prokaryot. a eukaryot. m.'
$ws.Range("B190").Value = 'ME301'
$ws.Range("A191").Value = 'Biologie půdních mikroorganismů'
$ws.Range("B191").Value = 'M101'
$ws.Range("A195").Value = 'Algoritmizace a programování II'
$ws.Range("B195").Value = 'APR2'
$ws.Range("A196").Value = 'Data Analysis and Visualisation'
$ws.Range("B196").Value = 'EDAV'
$ws.Range("A197").Value = 'Data Mining Techniq. based on R Software'
$ws.Range("B197").Value = 'EDMR'
$ws.Range("A198").Value = 'Algoritmizace a programování I'
$ws.Range("B198").Value = 'APR1'
$ws.Range("A199").Value = 'Python and R for Data Science'
$ws.Range("B199").Value = 'EPYR'
$ws.Range("A200").Value = 'Machine Learning Based on R Software'
$ws.Range("B200").Value = 'EMLR'
$ws.Range("A201").Value = 'Základy matematiky'
$ws.Range("B201").Value = 'P103'
$ws.Range("A202").Value = 'Pravděpodobnost a statistika I'
$ws.Range("B202").Value = 'K413'
$ws.Range("A203").Value = 'Vybrané partie z matematiky'
$ws.Range("B203").Value = 'P232'
$ws.Range("B206").Value = 'OONV'
$ws.Range("B207").Value = 'KOONV'
